# TrialsSetup.xlsx update: Power Query refresh added a new trial row
# ("REDEFINE HF") to the Sheet1!Query1 table, growing it from A1:B13 to
# A1:B14, and the refresh re-stamped column A (Trial Name) with an
# explicit "General" number format on the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Write the new trial name into the newly-added row. Progress (col B)
#    is left blank for this trial, matching the refreshed source data.
$ws.Range("A14").Value = "REDEFINE HF"

# 2. Re-apply the "General" number format across the table's Trial Name
#    column data rows (A2:A14) -- this is what the query refresh does to
#    every data cell in that column, old and new alike.
$ws.Range("A2:A14").NumberFormat = "General"

# 3. Grow the query table (ListObject "Query1") to include the new row.
$lo = $ws.ListObjects.Item("Query1")
$lo.Resize($ws.Range("A1:B14"))

# 4. Extend the hidden ExternalData_1 defined name that tracks the query
#    table's backing range so it also covers the new row.
$name = $wb.Names.Item("ExternalData_1")
$name.RefersTo = "=Sheet1!`$A`$1:`$B`$14"
